$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new allocations for "Micheal Rebourne" (row 5): Login, SendMessage, BanFromGroup
$ws.Range("B5").Value = "Login"
$ws.Range("C5").Value = "SendMessage"
$ws.Range("D5").Value = "BanFromGroup"

# Update selection to match the new active cell
$ws.Range("D5").Select()
